# Alterado Diagrama de Sequência Adicionar Pintor
# Swap the text of the two last "Flow of Events" steps (6 and 7) of the
# "Adicionar Pintor" use-case table, matching the updated sequence diagram,
# and update the lowercase wording of the "Devolve dados de autenticação do
# pintor" step. Then leave the active selection on D15, as last edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 6 (row 13) now reads what used to be step 7's text.
$ws.Range("D13").Value = "Adiciona Pintor a Lista de Pintores  "

# Step 7 (row 14) now reads the (re-cased) former step 6 text.
$ws.Range("D14").Value = "Devolve dados de autenticação do pintor"

# Match the saved selection state (active cell D15).
$ws.Range("D15").Select()
